$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume (E) columns are treated as plain text so that
# values such as "43.493.19" or "0.510" are not auto-converted to numbers/dates.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.493.19'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.38%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.333.27'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.54%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '305.44'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.33%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '101.75'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -2.29%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.510'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -2.86%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.509'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -2.21%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.27'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.37%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0797'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.92%  '

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.58%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.80'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.69%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.693.50'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.55%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.68'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.23%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.327.16'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.70%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.806'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.11%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.386.95'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.19%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.80'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.59%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0908'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.88%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.10'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.86%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.29'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.14%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.60'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.62%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.98'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -3.39%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.54'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -3.25%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.06%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.89'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.79%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.18'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -3.27%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '34.77'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -4.83%  '

$ws.Range("B30").Value = 'Cosmos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.23'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -3.73%  '

$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '164.79'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.64%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.03%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.07'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.86%  '

$ws.Range("B34").Value = 'RenderToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.62'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.69%  '

$ws.Range("B35").Value = 'WEMIXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.42'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -4.84%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.10'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -6.49%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0706'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -4.32%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.92'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -5.55%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.82'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -4.25%  '

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -3.84%  '

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -3.11%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.61'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +11.74%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.975.97'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.30%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0285'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.27%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '18.59'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -6.80%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.13'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.80%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.92'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -4.58%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '55.92'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -5.06%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.81'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.84%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.57'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.13%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.554.33'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.17%  '
